# adding remove card feature
# Populate the worksheet with the card data (row 1 and row 3, row 2 left blank)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 123
$ws.Range("B1").Value = 45

$ws.Range("A3").Value = 333
$ws.Range("B3").Value = 33

$ws.Range("H10").Select()
